# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" rows (16-21) are re-sorted into ascending period order
# (2108, 2109, 2110, 2111, 2112, 2201 instead of 2201, 2112, 2111, 2110,
# 2109, 2108), the "Valor Mora" (column F) figure that used to belong to
# period 2201 now travels with it to the bottom row (and vice-versa for the
# value that belonged to 2108), and every "Salario Basico" (column G) entry
# is refreshed from 908526 to 877803.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo Mora (column E) - now listed oldest -> newest
$ws.Range("E16").Value = "2108"
$ws.Range("E17").Value = "2109"
$ws.Range("E18").Value = "2110"
$ws.Range("E19").Value = "2111"
$ws.Range("E20").Value = "2112"
$ws.Range("E21").Value = "2201"

# Valor Mora (column F) - the two outer values swap rows, the rest stay put
$ws.Range("F16").Value = 35112
$ws.Range("F17").Value = 35112
$ws.Range("F18").Value = 35112
$ws.Range("F19").Value = 35112
$ws.Range("F20").Value = 35112
$ws.Range("F21").Value = 23408

# Salario Basico (column G) - updated for every period row
$ws.Range("G16").Value = 877803
$ws.Range("G17").Value = 877803
$ws.Range("G18").Value = 877803
$ws.Range("G19").Value = 877803
$ws.Range("G20").Value = 877803
$ws.Range("G21").Value = 877803
